# regenerace dokumentace 3. iterace
# - s ohledem na opravu generování dokumentace

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new line item as the next row after the last used row. Column B
# already carries its number format at the column level, so the new cell
# picks it up without any extra formatting step.
$ws.Range("A46").Value = "regenerace dokumentace 3. iterace"
$ws.Range("B46").Value = 0.5

# Move the frozen-pane viewport and active selection down to follow the
# newly-added row, same as Excel does when you scroll/select further down.
$excel.ActiveWindow.ScrollRow = 34
$ws.Range("A47").Select()
